$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.707.69'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -3.54%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.741.60'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -5.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.91'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -8.24%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4911'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -6.98%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.82'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -7.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2415'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -23.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05997'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -11.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.739.41'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -5.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06793'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -12.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.81'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -21.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.445'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -11.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.61'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -13.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.5822'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -25.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.0000'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.730.44'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.49'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -17.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006394'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -19.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.960.88'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -5.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.941'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -14.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.096'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -14.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.847'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -15.82%  '
$ws.Range("E26").Value = '  -4.52%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.848'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -16.72%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.457'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -13.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.47'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -15.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '99.94'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -9.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08101'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -7.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.726'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -11.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.374'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -17.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04365'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -10.63%  '
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.695'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.020'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -10.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6022'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -17.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.732'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -11.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.054'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -9.96%  '
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.15'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -6.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01485'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -14.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7858'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -12.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3805'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -20.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.136'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -13.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.012'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -21.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05092'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -12.35%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.24'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -13.09%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1065'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -14.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.26'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -12.77%  '
